# city prefer, settlement prefer finsihed, harbor prefer in process
#
# Catan "point_point" init data: column A (Point1) is being renumbered for
# rows 8-11 (each value decremented by one), and the sheet's last-saved
# selection moves from the old C2:D73 block down to a single cell, E4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 8-11: A column goes 5,6,7,8 -> 4,5,6,7 (B column is untouched)
$ws.Cells.Item(8, 1).Value = 4
$ws.Cells.Item(9, 1).Value = 5
$ws.Cells.Item(10, 1).Value = 6
$ws.Cells.Item(11, 1).Value = 7

# Move the active selection to E4 (was C2:D73 with active cell C2)
$ws.Range("E4").Select()

# Best-effort: also mirror the saved window geometry from the workbook view
# (xWindow/yWindow/windowWidth/windowHeight in the xml) onto the live window.
$excel.ActiveWindow.Left = 14400
$excel.ActiveWindow.Top = 0
$excel.ActiveWindow.Width = 14400
$excel.ActiveWindow.Height = 15600
